# Update tab names in all BOMs, fix bi-color LED naming.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab from "LOGx" to "BOM".
$ws.Name = "BOM"

# The "Ref" column cells below used an explicit (but visually identical)
# cell style; drop the redundant style override so the cells fall back to
# the default "Normal" style, matching the cleaned-up styles table.
$refCells = @("D3", "D4", "D6", "D10", "D17", "D18", "D20")
foreach ($addr in $refCells) {
    $ws.Range($addr).Style = "Normal"
}
